# Adapt column header formatting to respective input file names.
# "_old" -> "_FV2210", "_new" -> "_FV2304" for the 10 AHB attribute columns,
# then wrap the used range in a table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J carry the "_old" -> "_FV2210" headers, columns L-U carry the
# "_new" -> "_FV2304" headers. Column K ("diff") is left untouched.
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $leftCol = $i + 1
    $rightCol = $i + 12

    $ws.Cells.Item(1, $leftCol).Value = $oldHeaders[$i] + "_FV2210"
    $ws.Cells.Item(1, $rightCol).Value = $oldHeaders[$i] + "_FV2304"
}

# Turn the data range into an Excel table ("Table1") with the renamed headers.
$dataRange = $ws.Range("A1:U69")
$tbl = $ws.ListObjects.Add(1, $dataRange, $false, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, top-left of scrollable area A2).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Headers renamed, table created, header row frozen."
